$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Replace the "FINALIDADE" paragraph text.
# ------------------------------------------------------------------
$old1 = "Estabelecer processos, conforme DCA 16-5, que dependem de atividades da Subdivisão de Qualificação (NQUA) e setores relacionados, de acordo com a competência estabelecida pelo Regimento Interno do CELOG (RICA 21-34)."
$new1 = "Relacionar processos que compreendam as competências estabelecidas pelo Regimento Interno do CELOG (RICA 21-34) para a Subdivisão de Qualificação (NQUA), de forma a atender ao preconizado na DCA 16-5 – Gestão por processos no COMAER."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ------------------------------------------------------------------
# 2) Remove the extra blank (justified) paragraph that sits right
#    before the "PROCESSOS RELACIONADOS" heading, and rename that
#    heading to "ATRIBUIÇÕES".
# ------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "PROCESSOS RELACIONADOS") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 1) {
    $prevPara = $d.Paragraphs.Item($targetIndex - 1)
    if ($prevPara.Range.Text.Trim() -eq "") {
        $prevPara.Range.Delete()
    }
}

$d.Content.Find.Execute("PROCESSOS RELACIONADOS", $true, $false, $false, $false, $false, $true, 1, $false, "ATRIBUIÇÕES", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Update the first related-process description.
# ------------------------------------------------------------------
$old3 = "Gestão de certificação de material nacionalizado junto ao IFI"
$new3 = "Gestão de certificação/qualificação de material nacionalizado"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# ------------------------------------------------------------------
# 4) Renumber PLOG0007 to PLOG0010.
# ------------------------------------------------------------------
$d.Content.Find.Execute("PLOG0007", $true, $false, $false, $false, $false, $true, 1, $false, "PLOG0010", 2) | Out-Null

# ------------------------------------------------------------------
# 5) Insert a new related-process bullet right after PLOG0010's
#    paragraph, reusing its paragraph/list formatting.
# ------------------------------------------------------------------
$plog0010Index = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.IndexOf("PLOG0010") -ge 0) {
        $plog0010Index = $i
        break
    }
}

if ($plog0010Index -gt 0) {
    $plog0010Para = $d.Paragraphs.Item($plog0010Index)
    $plog0010Para.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($plog0010Index + 1)
    $newPara.Range.InsertBefore("PLOG0011 - Homologação de ensaio de controle de qualidade de fornecedor")
}
